$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$ws.Range("A169").Value = 46047.96788326389
$ws.Range("B169").Value = "AU001YBE"
$ws.Range("C169").Value = "Yobe"
$ws.Range("D169").Value = "Potiskum"
$ws.Range("E169").Value = "Sorghum"
$ws.Range("F169").Value = 21000
$ws.Range("G169").Value = 95
$ws.Range("H169").Value = 221.052631578947
$ws.Range("I169").Value = "high"
$ws.Range("J169").Value = "New"
$ws.Range("K169").Value = 224.052631578947

$ws.Range("A170").Value = 46047.96671715278
$ws.Range("B170").Value = "AU001YBE"
$ws.Range("C170").Value = "Yobe"
$ws.Range("D170").Value = "Potiskum"
$ws.Range("E170").Value = "Rice Processed"
$ws.Range("F170").Value = 92000
$ws.Range("G170").Value = 115
$ws.Range("H170").Value = 800
$ws.Range("I170").Value = "high"
$ws.Range("J170").Value = "New"
$ws.Range("K170").Value = 803

$ws.Range("A171").Value = 46047.96462563657
$ws.Range("B171").Value = "AU001YBE"
$ws.Range("C171").Value = "Yobe"
$ws.Range("D171").Value = "Potiskum"
$ws.Range("E171").Value = "Millet"
$ws.Range("F171").Value = 31000
$ws.Range("G171").Value = 105
$ws.Range("H171").Value = 295.23809523809501
$ws.Range("I171").Value = "high"
$ws.Range("J171").Value = "New"
$ws.Range("K171").Value = 298.23809523809501

$ws.Range("A172").Value = 46047.95940414352
$ws.Range("B172").Value = "AU001YBE"
$ws.Range("C172").Value = "Yobe"
$ws.Range("D172").Value = "Potiskum"
$ws.Range("E172").Value = "Maize"
$ws.Range("F172").Value = 28000
$ws.Range("G172").Value = 100
$ws.Range("H172").Value = 280
$ws.Range("I172").Value = "high"
$ws.Range("J172").Value = "New"
$ws.Range("K172").Value = 283

$ws.Range("A173").Value = 46047.95848561342
$ws.Range("B173").Value = "AU001YBE"
$ws.Range("C173").Value = "Yobe"
$ws.Range("D173").Value = "Potiskum"
$ws.Range("E173").Value = "Groundnut Gargaja"
$ws.Range("F173").Value = 110000
$ws.Range("G173").Value = 98
$ws.Range("H173").Value = 1122.44897959183
$ws.Range("I173").Value = "high"
$ws.Range("J173").Value = "New"
$ws.Range("K173").Value = 1125.44897959183

$ws.Range("A174").Value = 46047.9580441088
$ws.Range("B174").Value = "AU001YBE"
$ws.Range("C174").Value = "Yobe"
$ws.Range("D174").Value = "Potiskum"
$ws.Range("E174").Value = "Cowpea White"
$ws.Range("F174").Value = 50000
$ws.Range("G174").Value = 100
$ws.Range("H174").Value = 500
$ws.Range("I174").Value = "high"
$ws.Range("J174").Value = "New"
$ws.Range("K174").Value = 503

$ws.Range("A175").Value = 46046.72683696759
$ws.Range("B175").Value = "MH001GMB"
$ws.Range("C175").Value = "Gombe"
$ws.Range("D175").Value = "Biliri"
$ws.Range("E175").Value = "Soya Beans"
$ws.Range("F175").Value = 58000
$ws.Range("G175").Value = 96
$ws.Range("H175").Value = 604.16666666666595
$ws.Range("I175").Value = "medium"
$ws.Range("J175").Value = "New"
$ws.Range("K175").Value = 607.16666666666595

$ws.Range("A176").Value = 46046.72597743056
$ws.Range("B176").Value = "MH001GMB"
$ws.Range("C176").Value = "Gombe"
$ws.Range("D176").Value = "Biliri"
$ws.Range("E176").Value = "Sorghum Red"
$ws.Range("F176").Value = 28000
$ws.Range("G176").Value = 100
$ws.Range("H176").Value = 280
$ws.Range("I176").Value = "low"
$ws.Range("J176").Value = "New"
$ws.Range("K176").Value = 283

$ws.Range("A177").Value = 46046.72463974537
$ws.Range("B177").Value = "MH001GMB"
$ws.Range("C177").Value = "Gombe"
$ws.Range("D177").Value = "Biliri"
$ws.Range("E177").Value = "Rice Paddy"
$ws.Range("F177").Value = 27000
$ws.Range("G177").Value = 70
$ws.Range("H177").Value = 385.71428571428498
$ws.Range("I177").Value = "medium"
$ws.Range("J177").Value = "New"
$ws.Range("K177").Value = 388.71428571428498

$ws.Range("A178").Value = 46046.72378166667
$ws.Range("B178").Value = "MH001GMB"
$ws.Range("C178").Value = "Gombe"
$ws.Range("D178").Value = "Biliri"
$ws.Range("E178").Value = "Millet"
$ws.Range("F178").Value = 29000
$ws.Range("G178").Value = 100
$ws.Range("H178").Value = 290
$ws.Range("I178").Value = "low"
$ws.Range("J178").Value = "New"
$ws.Range("K178").Value = 293

$ws.Range("A179").Value = 46046.72292863426
$ws.Range("B179").Value = "MH001GMB"
$ws.Range("C179").Value = "Gombe"
$ws.Range("D179").Value = "Biliri"
$ws.Range("E179").Value = "Maize White"
$ws.Range("F179").Value = 22000
$ws.Range("G179").Value = 95
$ws.Range("H179").Value = 231.57894736842101
$ws.Range("I179").Value = "medium"
$ws.Range("J179").Value = "New"
$ws.Range("K179").Value = 234.57894736842101

$ws.Range("A180").Value = 46046.72207498843
$ws.Range("B180").Value = "MH001GMB"
$ws.Range("C180").Value = "Gombe"
$ws.Range("D180").Value = "Biliri"
$ws.Range("E180").Value = "Groundnut Kampala"
$ws.Range("F180").Value = 90000
$ws.Range("G180").Value = 85
$ws.Range("H180").Value = 1058.8235294117601
$ws.Range("I180").Value = "low"
$ws.Range("J180").Value = "New"
$ws.Range("K180").Value = 1061.8235294117601

$ws.Range("A181").Value = 46046.72028496527
$ws.Range("B181").Value = "MH001GMB"
$ws.Range("C181").Value = "Gombe"
$ws.Range("D181").Value = "Biliri"
$ws.Range("E181").Value = "Groundnut Gargaja"
$ws.Range("F181").Value = 85000
$ws.Range("G181").Value = 85
$ws.Range("H181").Value = 1000
$ws.Range("I181").Value = "medium"
$ws.Range("J181").Value = "New"
$ws.Range("K181").Value = 1003

$ws.Range("A182").Value = 46046.71930752315
$ws.Range("B182").Value = "MH001GMB"
$ws.Range("C182").Value = "Gombe"
$ws.Range("D182").Value = "Biliri"
$ws.Range("E182").Value = "Cowpea White"
$ws.Range("F182").Value = 43000
$ws.Range("G182").Value = 80
$ws.Range("H182").Value = 537.5
$ws.Range("I182").Value = "medium"
$ws.Range("J182").Value = "New"
$ws.Range("K182").Value = 540.5

$ws.Range("A183").Value = 46046.49286483796
$ws.Range("B183").Value = "IS001KDN"
$ws.Range("C183").Value = "Kaduna"
$ws.Range("D183").Value = "Pambegua"
$ws.Range("E183").Value = "Maize"
$ws.Range("F183").Value = 24000
$ws.Range("G183").Value = 237
$ws.Range("H183").Value = 101.26582278481
$ws.Range("I183").Value = "high"
$ws.Range("J183").Value = "New"
$ws.Range("K183").Value = 104.26582278481

$ws.Range("A184").Value = 46046.49192572917
$ws.Range("B184").Value = "IS001KDN"
$ws.Range("C184").Value = "Kaduna"
$ws.Range("D184").Value = "Pambegua"
$ws.Range("E184").Value = "Rice Paddy"
$ws.Range("F184").Value = 35000
$ws.Range("G184").Value = 380
$ws.Range("H184").Value = 92.105263157894697
$ws.Range("I184").Value = "medium"
$ws.Range("J184").Value = "New"
$ws.Range("K184").Value = 95.105263157894697

$ws.Range("A185").Value = 46046.4911015625
$ws.Range("B185").Value = "IS001KDN"
$ws.Range("C185").Value = "Kaduna"
$ws.Range("D185").Value = "Pambegua"
$ws.Range("E185").Value = "Soya Beans"
$ws.Range("F185").Value = 60000
$ws.Range("G185").Value = 630
$ws.Range("H185").Value = 95.238095238095198
$ws.Range("I185").Value = "medium"
$ws.Range("J185").Value = "New"
$ws.Range("K185").Value = 98.238095238095198

# Reset style (no explicit number format) for numeric columns F, G, H, K so they
# do not inherit the column-level style (matches freshly pasted/raw data cells).
$ws.Range("F169:H185").Style = "Normal"
$ws.Range("K169:K185").Style = "Normal"

# Column A holds timestamps formatted as dates - restore style/format after the reset above
$ws.Range("A169:A185").Style = "Normal"
$ws.Range("A169:A185").NumberFormat = "yyyy\-mm\-dd"

# Update the frozen-pane view / selection to reflect scrolling down to the newly added rows
$ws.Range("H164").Select()
